# Generate Report for Handoff
# Update localization status from "In Translation" to "Ready for handoff"
# and refresh the related timestamps, across the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: zh-cn / de-de status columns + generate-date column ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Ready for handoff"
$wsOverview.Range("F2").Value = "Ready for handoff"
$wsOverview.Range("G2").Value = "2016-08-19 22:45:51"

# --- zh-cn sheet: Status + Latest Handoff Datetime ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = "Ready for handoff"
$wsZh.Range("H2").Value = "2016-08-19 22:45:48"

# --- de-de sheet: Status + Latest Handoff Datetime ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = "Ready for handoff"
$wsDe.Range("H2").Value = "2016-08-19 22:45:51"

# --- Column width adjustments (columns widened to fit the longer status text) ---
$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33
$wsZh.Columns.Item(3).ColumnWidth = 16.33
$wsDe.Columns.Item(3).ColumnWidth = 16.33
